$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove column M (the data column holds the "white wine" style middle counter);
# column N shifts left into M, matching "remove column from alcohol data".
$ws.Columns.Item(13).Delete()

# Keep the active selection on the now-last data column (M1), matching the
# post-edit workbook state.
$ws.Range("M1").Select()
